$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Cell = "H138"; Value = 2587.1013 }
    @{ Sheet = "ALC"; Cell = "I138"; Value = 2146.8445 }
    @{ Sheet = "ALC"; Cell = "J138"; Value = 3169.7942 }
    @{ Sheet = "ALC"; Cell = "K138"; Value = 6440.5335 }
    @{ Sheet = "ALC"; Cell = "L138"; Value = 9509.382599999999 }
    @{ Sheet = "ALC"; Cell = "M138"; Value = -1300.5335 }
    @{ Sheet = "ALC"; Cell = "N138"; Value = -19789.3826 }
    @{ Sheet = "ARM"; Cell = "H74"; Value = 12449.154 }
    @{ Sheet = "ARM"; Cell = "I74"; Value = 1391.2858 }
    @{ Sheet = "ARM"; Cell = "K74"; Value = 1391.2858 }
    @{ Sheet = "ARM"; Cell = "M74"; Value = -517.2858000000001 }
    @{ Sheet = "ARM"; Cell = "H77"; Value = 12449.154 }
    @{ Sheet = "ARM"; Cell = "I77"; Value = 1391.2858 }
    @{ Sheet = "ARM"; Cell = "K77"; Value = 6956.429 }
    @{ Sheet = "ARM"; Cell = "M77"; Value = -2588.429 }
    @{ Sheet = "ARM"; Cell = "H132"; Value = 23168000 }
    @{ Sheet = "ARM"; Cell = "I132"; Value = 43647756 }
    @{ Sheet = "ARM"; Cell = "J132"; Value = 128274 }
    @{ Sheet = "ARM"; Cell = "K132"; Value = 130943268 }
    @{ Sheet = "ARM"; Cell = "L132"; Value = 384822 }
    @{ Sheet = "ARM"; Cell = "M132"; Value = -130940738 }
    @{ Sheet = "ARM"; Cell = "N132"; Value = -389882 }
    @{ Sheet = "BSM"; Cell = "H105"; Value = 1648 }
    @{ Sheet = "BSM"; Cell = "I105"; Value = 1425.7142 }
    @{ Sheet = "BSM"; Cell = "J105"; Value = 2166.6667 }
    @{ Sheet = "BSM"; Cell = "K105"; Value = 1425.7142 }
    @{ Sheet = "BSM"; Cell = "L105"; Value = 2166.6667 }
    @{ Sheet = "BSM"; Cell = "M105"; Value = 321.2858000000001 }
    @{ Sheet = "BSM"; Cell = "N105"; Value = -5660.6667 }
    @{ Sheet = "BSM"; Cell = "H132"; Value = 60780 }
    @{ Sheet = "BSM"; Cell = "J132"; Value = 60780 }
    @{ Sheet = "BSM"; Cell = "L132"; Value = 60780 }
    @{ Sheet = "BSM"; Cell = "N132"; Value = -70900 }
    @{ Sheet = "CRP"; Cell = "H31"; Value = 7939839.5 }
    @{ Sheet = "CRP"; Cell = "J31"; Value = 16671985 }
    @{ Sheet = "CRP"; Cell = "L31"; Value = 16671985 }
    @{ Sheet = "CRP"; Cell = "N31"; Value = -16672575 }
    @{ Sheet = "CRP"; Cell = "H33"; Value = 19866.666 }
    @{ Sheet = "CRP"; Cell = "I33"; Value = 14950 }
    @{ Sheet = "CRP"; Cell = "J33"; Value = 29700 }
    @{ Sheet = "CRP"; Cell = "K33"; Value = 14950 }
    @{ Sheet = "CRP"; Cell = "L33"; Value = 29700 }
    @{ Sheet = "CRP"; Cell = "M33"; Value = -14571 }
    @{ Sheet = "CRP"; Cell = "N33"; Value = -30458 }
    @{ Sheet = "CRP"; Cell = "H34"; Value = 7939839.5 }
    @{ Sheet = "CRP"; Cell = "J34"; Value = 16671985 }
    @{ Sheet = "CRP"; Cell = "L34"; Value = 16671985 }
    @{ Sheet = "CRP"; Cell = "N34"; Value = -16672389 }
    @{ Sheet = "CRP"; Cell = "H41"; Value = 13040 }
    @{ Sheet = "CRP"; Cell = "I41"; Value = 2600 }
    @{ Sheet = "CRP"; Cell = "J41"; Value = 20000 }
    @{ Sheet = "CRP"; Cell = "K41"; Value = 2600 }
    @{ Sheet = "CRP"; Cell = "L41"; Value = 20000 }
    @{ Sheet = "CRP"; Cell = "M41"; Value = -2172 }
    @{ Sheet = "CRP"; Cell = "N41"; Value = -20856 }
    @{ Sheet = "CRP"; Cell = "H50"; Value = 10082.286 }
    @{ Sheet = "CRP"; Cell = "J50"; Value = 10082.286 }
    @{ Sheet = "CRP"; Cell = "L50"; Value = 10082.286 }
    @{ Sheet = "CRP"; Cell = "N50"; Value = -11332.286 }
    @{ Sheet = "CRP"; Cell = "H51"; Value = 9490.454 }
    @{ Sheet = "CRP"; Cell = "J51"; Value = 9529.5 }
    @{ Sheet = "CRP"; Cell = "L51"; Value = 9529.5 }
    @{ Sheet = "CRP"; Cell = "N51"; Value = -11001.5 }
    @{ Sheet = "CRP"; Cell = "H55"; Value = 950 }
    @{ Sheet = "CRP"; Cell = "I55"; Value = 950 }
    @{ Sheet = "CRP"; Cell = "K55"; Value = 950 }
    @{ Sheet = "CRP"; Cell = "M55"; Value = -635 }
    @{ Sheet = "CRP"; Cell = "H58"; Value = 2928.5605 }
    @{ Sheet = "CRP"; Cell = "I58"; Value = 1211.95 }
    @{ Sheet = "CRP"; Cell = "J58"; Value = 3674.913 }
    @{ Sheet = "CRP"; Cell = "K58"; Value = 1211.95 }
    @{ Sheet = "CRP"; Cell = "L58"; Value = 3674.913 }
    @{ Sheet = "CRP"; Cell = "M58"; Value = -1008.95 }
    @{ Sheet = "CRP"; Cell = "N58"; Value = -4080.913 }
    @{ Sheet = "CRP"; Cell = "H59"; Value = 15302.6 }
    @{ Sheet = "CRP"; Cell = "J59"; Value = 15225.111 }
    @{ Sheet = "CRP"; Cell = "L59"; Value = 15225.111 }
    @{ Sheet = "CRP"; Cell = "N59"; Value = -17515.111 }
    @{ Sheet = "CRP"; Cell = "H60"; Value = 10023.556 }
    @{ Sheet = "CRP"; Cell = "J60"; Value = 10023.556 }
    @{ Sheet = "CRP"; Cell = "L60"; Value = 10023.556 }
    @{ Sheet = "CRP"; Cell = "N60"; Value = -11045.556 }
    @{ Sheet = "CRP"; Cell = "H61"; Value = 9490.454 }
    @{ Sheet = "CRP"; Cell = "J61"; Value = 9529.5 }
    @{ Sheet = "CRP"; Cell = "L61"; Value = 9529.5 }
    @{ Sheet = "CRP"; Cell = "N61"; Value = -10225.5 }
    @{ Sheet = "CRP"; Cell = "H68"; Value = 17792 }
    @{ Sheet = "CRP"; Cell = "J68"; Value = 18183.555 }
    @{ Sheet = "CRP"; Cell = "L68"; Value = 18183.555 }
    @{ Sheet = "CRP"; Cell = "N68"; Value = -19681.555 }
    @{ Sheet = "CRP"; Cell = "H71"; Value = 17792 }
    @{ Sheet = "CRP"; Cell = "J71"; Value = 18183.555 }
    @{ Sheet = "CRP"; Cell = "L71"; Value = 54550.665 }
    @{ Sheet = "CRP"; Cell = "N71"; Value = -62038.665 }
    @{ Sheet = "CRP"; Cell = "H74"; Value = 14431.692 }
    @{ Sheet = "CRP"; Cell = "I74"; Value = 5085 }
    @{ Sheet = "CRP"; Cell = "J74"; Value = 16131.091 }
    @{ Sheet = "CRP"; Cell = "K74"; Value = 5085 }
    @{ Sheet = "CRP"; Cell = "L74"; Value = 16131.091 }
    @{ Sheet = "CRP"; Cell = "M74"; Value = -4211 }
    @{ Sheet = "CRP"; Cell = "N74"; Value = -17879.091 }
    @{ Sheet = "CRP"; Cell = "H77"; Value = 14431.692 }
    @{ Sheet = "CRP"; Cell = "I77"; Value = 5085 }
    @{ Sheet = "CRP"; Cell = "J77"; Value = 16131.091 }
    @{ Sheet = "CRP"; Cell = "K77"; Value = 15255 }
    @{ Sheet = "CRP"; Cell = "L77"; Value = 48393.273 }
    @{ Sheet = "CRP"; Cell = "M77"; Value = -10887 }
    @{ Sheet = "CRP"; Cell = "N77"; Value = -57129.273 }
    @{ Sheet = "CRP"; Cell = "H99"; Value = 144587.72 }
    @{ Sheet = "CRP"; Cell = "I99"; Value = 1934 }
    @{ Sheet = "CRP"; Cell = "J99"; Value = 251578 }
    @{ Sheet = "CRP"; Cell = "K99"; Value = 1934 }
    @{ Sheet = "CRP"; Cell = "L99"; Value = 251578 }
    @{ Sheet = "CRP"; Cell = "M99"; Value = -436 }
    @{ Sheet = "CRP"; Cell = "N99"; Value = -254574 }
    @{ Sheet = "CRP"; Cell = "H126"; Value = 144587.72 }
    @{ Sheet = "CRP"; Cell = "I126"; Value = 1934 }
    @{ Sheet = "CRP"; Cell = "J126"; Value = 251578 }
    @{ Sheet = "CRP"; Cell = "K126"; Value = 5802 }
    @{ Sheet = "CRP"; Cell = "L126"; Value = 754734 }
    @{ Sheet = "CRP"; Cell = "M126"; Value = -3332 }
    @{ Sheet = "CRP"; Cell = "N126"; Value = -759674 }
    @{ Sheet = "CRP"; Cell = "H136"; Value = 2928.5605 }
    @{ Sheet = "CRP"; Cell = "I136"; Value = 1211.95 }
    @{ Sheet = "CRP"; Cell = "J136"; Value = 3674.913 }
    @{ Sheet = "CRP"; Cell = "K136"; Value = 3635.85 }
    @{ Sheet = "CRP"; Cell = "L136"; Value = 11024.739 }
    @{ Sheet = "CRP"; Cell = "M136"; Value = -1085.85 }
    @{ Sheet = "CRP"; Cell = "N136"; Value = -16124.739 }
    @{ Sheet = "CUL"; Cell = "H37"; Value = 31008.637 }
    @{ Sheet = "CUL"; Cell = "J37"; Value = 31008.637 }
    @{ Sheet = "CUL"; Cell = "L37"; Value = 93025.91099999999 }
    @{ Sheet = "CUL"; Cell = "N37"; Value = -93249.91099999999 }
    @{ Sheet = "CUL"; Cell = "H49"; Value = 1100 }
    @{ Sheet = "CUL"; Cell = "J49"; Value = 1175 }
    @{ Sheet = "CUL"; Cell = "L49"; Value = 3525 }
    @{ Sheet = "CUL"; Cell = "N49"; Value = -3837 }
    @{ Sheet = "CUL"; Cell = "H109"; Value = 2264.7 }
    @{ Sheet = "CUL"; Cell = "I109"; Value = 1133.1666 }
    @{ Sheet = "CUL"; Cell = "J109"; Value = 2749.6428 }
    @{ Sheet = "CUL"; Cell = "K109"; Value = 3399.4998 }
    @{ Sheet = "CUL"; Cell = "L109"; Value = 8248.928400000001 }
    @{ Sheet = "CUL"; Cell = "M109"; Value = -2359.4998 }
    @{ Sheet = "CUL"; Cell = "N109"; Value = -10328.9284 }
    @{ Sheet = "CUL"; Cell = "H113"; Value = 548.1622 }
    @{ Sheet = "CUL"; Cell = "I113"; Value = 520.381 }
    @{ Sheet = "CUL"; Cell = "J113"; Value = 584.625 }
    @{ Sheet = "CUL"; Cell = "K113"; Value = 1561.143 }
    @{ Sheet = "CUL"; Cell = "L113"; Value = 1753.875 }
    @{ Sheet = "CUL"; Cell = "M113"; Value = 608.857 }
    @{ Sheet = "CUL"; Cell = "N113"; Value = -6093.875 }
    @{ Sheet = "GSM"; Cell = "H132"; Value = 66669570 }
    @{ Sheet = "GSM"; Cell = "I132"; Value = 125002024 }
    @{ Sheet = "GSM"; Cell = "J132"; Value = 3903.2856 }
    @{ Sheet = "GSM"; Cell = "K132"; Value = 375006072 }
    @{ Sheet = "GSM"; Cell = "L132"; Value = 11709.8568 }
    @{ Sheet = "GSM"; Cell = "M132"; Value = -375003542 }
    @{ Sheet = "GSM"; Cell = "N132"; Value = -16769.8568 }
    @{ Sheet = "LTW"; Cell = "H22"; Value = 366.16666 }
    @{ Sheet = "LTW"; Cell = "I22"; Value = 250.5 }
    @{ Sheet = "LTW"; Cell = "J22"; Value = 424 }
    @{ Sheet = "LTW"; Cell = "K22"; Value = 250.5 }
    @{ Sheet = "LTW"; Cell = "L22"; Value = 424 }
    @{ Sheet = "LTW"; Cell = "M22"; Value = 44.5 }
    @{ Sheet = "LTW"; Cell = "N22"; Value = -1014 }
    @{ Sheet = "LTW"; Cell = "H27"; Value = 366.16666 }
    @{ Sheet = "LTW"; Cell = "I27"; Value = 250.5 }
    @{ Sheet = "LTW"; Cell = "J27"; Value = 424 }
    @{ Sheet = "LTW"; Cell = "K27"; Value = 250.5 }
    @{ Sheet = "LTW"; Cell = "L27"; Value = 424 }
    @{ Sheet = "LTW"; Cell = "M27"; Value = -143.5 }
    @{ Sheet = "LTW"; Cell = "N27"; Value = -638 }
    @{ Sheet = "LTW"; Cell = "H43"; Value = 4000 }
    @{ Sheet = "LTW"; Cell = "I43"; Value = 4000 }
    @{ Sheet = "LTW"; Cell = "J43"; Value = 4000 }
    @{ Sheet = "LTW"; Cell = "K43"; Value = 4000 }
    @{ Sheet = "LTW"; Cell = "L43"; Value = 4000 }
    @{ Sheet = "LTW"; Cell = "M43"; Value = -3807 }
    @{ Sheet = "LTW"; Cell = "N43"; Value = -4386 }
    @{ Sheet = "LTW"; Cell = "H93"; Value = 1167.1666 }
    @{ Sheet = "LTW"; Cell = "I93"; Value = 827.8182 }
    @{ Sheet = "LTW"; Cell = "K93"; Value = 827.8182 }
    @{ Sheet = "LTW"; Cell = "M93"; Value = 420.1818 }
    @{ Sheet = "LTW"; Cell = "H127"; Value = 38543 }
    @{ Sheet = "LTW"; Cell = "J127"; Value = 38543 }
    @{ Sheet = "LTW"; Cell = "L127"; Value = 38543 }
    @{ Sheet = "LTW"; Cell = "N127"; Value = -48463 }
    @{ Sheet = "LTW"; Cell = "H132"; Value = 8930.842000000001 }
    @{ Sheet = "LTW"; Cell = "I132"; Value = 10279.2 }
    @{ Sheet = "LTW"; Cell = "J132"; Value = 3874.5 }
    @{ Sheet = "LTW"; Cell = "K132"; Value = 30837.6 }
    @{ Sheet = "LTW"; Cell = "L132"; Value = 11623.5 }
    @{ Sheet = "LTW"; Cell = "M132"; Value = -28307.6 }
    @{ Sheet = "LTW"; Cell = "N132"; Value = -16683.5 }
    @{ Sheet = "WVR"; Cell = "H100"; Value = 433.33334 }
    @{ Sheet = "WVR"; Cell = "I100"; Value = 362.5 }
    @{ Sheet = "WVR"; Cell = "J100"; Value = 1000 }
    @{ Sheet = "WVR"; Cell = "K100"; Value = 725 }
    @{ Sheet = "WVR"; Cell = "L100"; Value = 2000 }
    @{ Sheet = "WVR"; Cell = "M100"; Value = -184 }
    @{ Sheet = "WVR"; Cell = "N100"; Value = -3082 }
    @{ Sheet = "WVR"; Cell = "H107"; Value = 695.05 }
    @{ Sheet = "WVR"; Cell = "I107"; Value = 428.06668 }
    @{ Sheet = "WVR"; Cell = "J107"; Value = 1496 }
    @{ Sheet = "WVR"; Cell = "K107"; Value = 1284.20004 }
    @{ Sheet = "WVR"; Cell = "L107"; Value = 4488 }
    @{ Sheet = "WVR"; Cell = "M107"; Value = 635.7999599999998 }
    @{ Sheet = "WVR"; Cell = "N107"; Value = -8328 }
    @{ Sheet = "WVR"; Cell = "H132"; Value = 28024520 }
    @{ Sheet = "WVR"; Cell = "I132"; Value = 39233068 }
    @{ Sheet = "WVR"; Cell = "J132"; Value = 3151 }
    @{ Sheet = "WVR"; Cell = "K132"; Value = 117699204 }
    @{ Sheet = "WVR"; Cell = "L132"; Value = 9453 }
    @{ Sheet = "WVR"; Cell = "M132"; Value = -117696674 }
    @{ Sheet = "WVR"; Cell = "N132"; Value = -14513 }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $ws.Range($chg.Cell).Value = $chg.Value
}